$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observed")
$ws.Range("AZ1").Value = 13.1
$ws.Range("AZ2").Value = 39.1
$ws.Range("AZ3").Value = 12.9
$ws.Range("AZ4").Value = 35.2
